$d = $word.ActiveDocument

# --- 1) "2. Functional Requirements:" was split across 3 runs ("2", ".",
#        " Functional Requirements:") -- collapse back into a single run by
#        replacing the whole heading text (Find/Replace naturally emits one
#        fresh run for the matched range).
$rng = $d.Content
$rng.Find.Execute("2. Functional Requirements:", $false, $false, $false, $false, $false, $true, 1, $false, "2. Functional Requirements:", 2) | Out-Null

# --- 2) "3. ERD Diagram:" -> "3. ERM Diagram:", but now split across three
#        runs ("3. ER", "M", " Diagram:") all sharing identical formatting.
#        First fix the text itself with a normal Find/Replace ...
$rng2 = $d.Content
$rng2.Find.Execute("3. ERD Diagram:", $false, $false, $false, $false, $false, $true, 1, $false, "3. ERM Diagram:", 2) | Out-Null

# ... then re-find the resulting text and pin the two interior split points
# with (transient) bookmarks so the engine keeps the run boundaries instead
# of silently re-coalescing the identically-formatted neighbours.
$rng3 = $d.Content
$rng3.Find.Execute("3. ERM Diagram:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$headingStart = $rng3.Start

$splitAfterER = $headingStart + 5   # right after "3. ER"
$splitAfterM  = $headingStart + 6   # right after "3. ERM"

$bm1Range = $d.Range($splitAfterER, $splitAfterER)
$d.Bookmarks.Add("ztmpSplit1", $bm1Range) | Out-Null
$bm2Range = $d.Range($splitAfterM, $splitAfterM)
$d.Bookmarks.Add("ztmpSplit2", $bm2Range) | Out-Null

if ($d.Bookmarks.Exists("ztmpSplit1")) { $d.Bookmarks.Item("ztmpSplit1").Delete() }
if ($d.Bookmarks.Exists("ztmpSplit2")) { $d.Bookmarks.Item("ztmpSplit2").Delete() }

# --- 3) The ER diagram picture's run gains an explicit <w:rPr><w:noProof/></w:rPr>.
$shapes = $d.InlineShapes
if ($shapes.Count -ge 1) {
    $picRange = $shapes.Item(1).Range
    $picRange.NoProofing = 1
}

# --- 4) "5. Data Types for Attributes:" was split across 2 runs ("5.",
#        " Data Types for Attributes:") -- collapse back into one run.
$rng4 = $d.Content
$rng4.Find.Execute("5. Data Types for Attributes:", $false, $false, $false, $false, $false, $true, 1, $false, "5. Data Types for Attributes:", 2) | Out-Null

# --- 5) "6. Integrity Constraints:" was split across 2 runs ("6.",
#        " Integrity Constraints:") -- collapse back into one run.
$rng5 = $d.Content
$rng5.Find.Execute("6. Integrity Constraints:", $false, $false, $false, $false, $false, $true, 1, $false, "6. Integrity Constraints:", 2) | Out-Null

Write-Output "edit.ps1 complete"
